$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range('K2').Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474993058530/1743637847690/PaZ2GKc6aN-camera-video-segment-driver-1743637850190.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMK65G53T%2F20250403%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250403T180030Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIj%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCICgh8znk17t3TA%2BoKV6Gs6dr67pjsww%2B74TUfK2AKgf7AiEAz%2Bsh73bZQ2b5CCsKg0Vjdzko8tsT%2FHUnEee9To43qz0q5gMI8f%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDCtu7H5t94rCmyRPWSq6A40NpyKaxkLrPb21f0WvTGT%2BIGiyn7Vbahoust2%2FGVAklzVYhIFiLefEp7yFI%2B7oJscwdjq7YKspHYMijhn7XRWKB4q464bVv7AXR5Uk5%2FmcOBPpkQAl0Be5gWWOYqROqpt%2FfDBW4shLaQ%2B2Djfz8obOjliQROGSRIdhfFY2BupuWf1VQ8prpoAvXUduCB3rxfGPNwH89vgp9nigybAXOLa25dLV1oq%2BzBxb39qnxbzdW8VQMOeQCwvTX%2Fq0vDxuDeN0iQREY4vwy%2BL13iBqN%2BkPYF9Jctlv43thbAG6o58%2Fh%2BOQQtWOGAykxPLN4qYJNqUBr%2BXgd1UBqMrKKgyIrmtMyVDf1zZYCtD7TuGL2AIszRTuVQY5tEIt2kxWCKFfKSymj8xk6kI%2ByUwaDrztK2xDDGJ%2Bqpe7Lcn5m6MkSrqab%2BvdDQEk6vYmJNdixVjZT6Q7SzfLXAkUsHZGE76ZlcPr8ko7UkStEK06fM%2BCql0r%2BD%2FXGcC0jvV5hXcN3QACFNyeMargdFsewTxlm%2FpWJgfI7X4u0LO8goE%2F67DaL8MqGpWCVkbvMPrJVZTdQIEwip1Af0893leelaswhti6vwY6pQHL%2FVpWyuVSfKKmSWKrk0yqWCTU%2FIZQrKI6CNfDOMGVTVeh%2BeQlMGeneOMBX83X3%2FtDD9GAWK2Grz19JgJk9GJufgKV75gfSafxWlhJ2gBKoWQZYRuoK%2Fo8ietECWT19yHU6Xs%2BbtZFlDaa346Gp%2BkCKvs3RDD%2FQ%2Bs%2BDQYMcujj%2B3UiGHSGthFgBpaAKoLr05Qg0%2BoOLYJ3JHH1JZtQYDvcgCTG6BY%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2004%20Apr%202025%2002%3A00%3A30%20GMT&X-Amz-Signature=5a65584d656c7614d0dddf5d4a6107fc092badea9ec9da03ac23e4eb717f7159'
$ws.Range('K3').Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474993058530/1743634721601/5dyCDJLTC3-camera-video-segment-driver-1743634724101.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMK65G53T%2F20250403%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250403T180030Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIj%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCICgh8znk17t3TA%2BoKV6Gs6dr67pjsww%2B74TUfK2AKgf7AiEAz%2Bsh73bZQ2b5CCsKg0Vjdzko8tsT%2FHUnEee9To43qz0q5gMI8f%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDCtu7H5t94rCmyRPWSq6A40NpyKaxkLrPb21f0WvTGT%2BIGiyn7Vbahoust2%2FGVAklzVYhIFiLefEp7yFI%2B7oJscwdjq7YKspHYMijhn7XRWKB4q464bVv7AXR5Uk5%2FmcOBPpkQAl0Be5gWWOYqROqpt%2FfDBW4shLaQ%2B2Djfz8obOjliQROGSRIdhfFY2BupuWf1VQ8prpoAvXUduCB3rxfGPNwH89vgp9nigybAXOLa25dLV1oq%2BzBxb39qnxbzdW8VQMOeQCwvTX%2Fq0vDxuDeN0iQREY4vwy%2BL13iBqN%2BkPYF9Jctlv43thbAG6o58%2Fh%2BOQQtWOGAykxPLN4qYJNqUBr%2BXgd1UBqMrKKgyIrmtMyVDf1zZYCtD7TuGL2AIszRTuVQY5tEIt2kxWCKFfKSymj8xk6kI%2ByUwaDrztK2xDDGJ%2Bqpe7Lcn5m6MkSrqab%2BvdDQEk6vYmJNdixVjZT6Q7SzfLXAkUsHZGE76ZlcPr8ko7UkStEK06fM%2BCql0r%2BD%2FXGcC0jvV5hXcN3QACFNyeMargdFsewTxlm%2FpWJgfI7X4u0LO8goE%2F67DaL8MqGpWCVkbvMPrJVZTdQIEwip1Af0893leelaswhti6vwY6pQHL%2FVpWyuVSfKKmSWKrk0yqWCTU%2FIZQrKI6CNfDOMGVTVeh%2BeQlMGeneOMBX83X3%2FtDD9GAWK2Grz19JgJk9GJufgKV75gfSafxWlhJ2gBKoWQZYRuoK%2Fo8ietECWT19yHU6Xs%2BbtZFlDaa346Gp%2BkCKvs3RDD%2FQ%2Bs%2BDQYMcujj%2B3UiGHSGthFgBpaAKoLr05Qg0%2BoOLYJ3JHH1JZtQYDvcgCTG6BY%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2004%20Apr%202025%2002%3A00%3A30%20GMT&X-Amz-Signature=6e23ea3c892bf96931b3f9cd71844d871054bbb449517535acfa30c3681ab615'
$ws.Range('K5').Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991109864/1743618603072/lsfY7Elc8s-camera-video-segment-driver-1743618608072.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMK65G53T%2F20250403%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250403T180030Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIj%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCICgh8znk17t3TA%2BoKV6Gs6dr67pjsww%2B74TUfK2AKgf7AiEAz%2Bsh73bZQ2b5CCsKg0Vjdzko8tsT%2FHUnEee9To43qz0q5gMI8f%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDCtu7H5t94rCmyRPWSq6A40NpyKaxkLrPb21f0WvTGT%2BIGiyn7Vbahoust2%2FGVAklzVYhIFiLefEp7yFI%2B7oJscwdjq7YKspHYMijhn7XRWKB4q464bVv7AXR5Uk5%2FmcOBPpkQAl0Be5gWWOYqROqpt%2FfDBW4shLaQ%2B2Djfz8obOjliQROGSRIdhfFY2BupuWf1VQ8prpoAvXUduCB3rxfGPNwH89vgp9nigybAXOLa25dLV1oq%2BzBxb39qnxbzdW8VQMOeQCwvTX%2Fq0vDxuDeN0iQREY4vwy%2BL13iBqN%2BkPYF9Jctlv43thbAG6o58%2Fh%2BOQQtWOGAykxPLN4qYJNqUBr%2BXgd1UBqMrKKgyIrmtMyVDf1zZYCtD7TuGL2AIszRTuVQY5tEIt2kxWCKFfKSymj8xk6kI%2ByUwaDrztK2xDDGJ%2Bqpe7Lcn5m6MkSrqab%2BvdDQEk6vYmJNdixVjZT6Q7SzfLXAkUsHZGE76ZlcPr8ko7UkStEK06fM%2BCql0r%2BD%2FXGcC0jvV5hXcN3QACFNyeMargdFsewTxlm%2FpWJgfI7X4u0LO8goE%2F67DaL8MqGpWCVkbvMPrJVZTdQIEwip1Af0893leelaswhti6vwY6pQHL%2FVpWyuVSfKKmSWKrk0yqWCTU%2FIZQrKI6CNfDOMGVTVeh%2BeQlMGeneOMBX83X3%2FtDD9GAWK2Grz19JgJk9GJufgKV75gfSafxWlhJ2gBKoWQZYRuoK%2Fo8ietECWT19yHU6Xs%2BbtZFlDaa346Gp%2BkCKvs3RDD%2FQ%2Bs%2BDQYMcujj%2B3UiGHSGthFgBpaAKoLr05Qg0%2BoOLYJ3JHH1JZtQYDvcgCTG6BY%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2004%20Apr%202025%2002%3A00%3A30%20GMT&X-Amz-Signature=2e84f76a5dcdab43c82304b6bb80a10ea640a853ae9e3464d5a82ede5b73924c'
$ws.Range('L5').Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109864/1743618603072/krpw3PxSiC-camera-video-segment-1743618608072.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMK65G53T%2F20250403%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250403T180030Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIj%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCICgh8znk17t3TA%2BoKV6Gs6dr67pjsww%2B74TUfK2AKgf7AiEAz%2Bsh73bZQ2b5CCsKg0Vjdzko8tsT%2FHUnEee9To43qz0q5gMI8f%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDCtu7H5t94rCmyRPWSq6A40NpyKaxkLrPb21f0WvTGT%2BIGiyn7Vbahoust2%2FGVAklzVYhIFiLefEp7yFI%2B7oJscwdjq7YKspHYMijhn7XRWKB4q464bVv7AXR5Uk5%2FmcOBPpkQAl0Be5gWWOYqROqpt%2FfDBW4shLaQ%2B2Djfz8obOjliQROGSRIdhfFY2BupuWf1VQ8prpoAvXUduCB3rxfGPNwH89vgp9nigybAXOLa25dLV1oq%2BzBxb39qnxbzdW8VQMOeQCwvTX%2Fq0vDxuDeN0iQREY4vwy%2BL13iBqN%2BkPYF9Jctlv43thbAG6o58%2Fh%2BOQQtWOGAykxPLN4qYJNqUBr%2BXgd1UBqMrKKgyIrmtMyVDf1zZYCtD7TuGL2AIszRTuVQY5tEIt2kxWCKFfKSymj8xk6kI%2ByUwaDrztK2xDDGJ%2Bqpe7Lcn5m6MkSrqab%2BvdDQEk6vYmJNdixVjZT6Q7SzfLXAkUsHZGE76ZlcPr8ko7UkStEK06fM%2BCql0r%2BD%2FXGcC0jvV5hXcN3QACFNyeMargdFsewTxlm%2FpWJgfI7X4u0LO8goE%2F67DaL8MqGpWCVkbvMPrJVZTdQIEwip1Af0893leelaswhti6vwY6pQHL%2FVpWyuVSfKKmSWKrk0yqWCTU%2FIZQrKI6CNfDOMGVTVeh%2BeQlMGeneOMBX83X3%2FtDD9GAWK2Grz19JgJk9GJufgKV75gfSafxWlhJ2gBKoWQZYRuoK%2Fo8ietECWT19yHU6Xs%2BbtZFlDaa346Gp%2BkCKvs3RDD%2FQ%2Bs%2BDQYMcujj%2B3UiGHSGthFgBpaAKoLr05Qg0%2BoOLYJ3JHH1JZtQYDvcgCTG6BY%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2004%20Apr%202025%2002%3A00%3A30%20GMT&X-Amz-Signature=9f0fd146a7f161f9b3280bcead70c0acba458fc6926da7d70eb4c6f7f82d2eed'
